# The workbook's row 2 and row 3 contain two observation records that need
# to have their distinguishing field values swapped (columns A, B, D, E, F,
# G, H, I, J and AI). All other columns already hold identical values in
# both rows, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "D", "E", "F", "G", "H", "I", "J", "AI")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value()
    $val3 = $cell3.Value()

    if ($val2 -eq $null) { $val2 = "" }
    if ($val3 -eq $null) { $val3 = "" }

    $cell2.Value = $val3
    $cell3.Value = $val2
}
